$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Construction")

# --- Add a second week column ("18/7") next to the existing "17/7" ---
$ws.Range("F3").Value = "18/7"

# --- Week 18/7 "Plus"/"Minus" figures for each team member ---
# Duy (rows 4-5): Plus unchanged (blank), Minus 17/7=10, 18/7=5 (new)
$ws.Range("F5").Value = 5

# Truong (rows 6-7): Plus 18/7=2 (new), Minus 17/7=5, 18/7 = 3+2 (new formula)
$ws.Range("E6").Value = 2
$ws.Range("F7").Formula = "=3+2"

# Giang (rows 8-9): Plus unchanged (blank), Minus 17/7=5, 18/7 = 3+2 (new formula)
$ws.Range("F9").Formula = "=3+2"

# Hai (rows 10-11): Plus unchanged (blank), Minus 17/7=5, 18/7=3 (new)
$ws.Range("F11").Value = 3

# --- Comments on the reviewed cells ---
$ws.Range("E6").AddComment("Have good question to make clear issue about Portlet 1.0, 2.0") | Out-Null
$ws.Range("F7").AddComment("-2: Commit unneccessary folder into SVN “http://oopms.googlecode.com/svn/trunk/SourceCode/RequirementModule/build”") | Out-Null
$ws.Range("F9").AddComment("-2: Commit unneccessary folder into SVN “http://oopms.googlecode.com/svn/trunk/SourceCode/PlannerModule/build`n-3: File .mpp has not updated") | Out-Null

# --- Notes section below the table ---
# "Note" label moves from D13 to A13 and becomes "Note:", bold
$ws.Range("D13").ClearContents()
$ws.Range("A13").Value = "Note:"
$ws.Range("A13").Font.Bold = $true

# Existing note stays, new note added below it
$ws.Range("E14").Value = "No update plan weekly"
$ws.Range("F15").Value = "Minus: File .mpp has not updated"

# New "Pending issues:" section, bold heading + red bullet points
$ws.Range("A17").Value = "Pending issues:"
$ws.Range("A17").Font.Bold = $true

$ws.Range("B18").Value = "Commit unneccessary into the SVN such as folder “build” within Project"
$ws.Range("B18").Font.Color = 255

$ws.Range("B19").Value = "The file .mpp has not updated % Completion of tasks"
$ws.Range("B19").Font.Color = 255

# --- Selection moves to A14 ---
$ws.Activate()
[void]$ws.Range("A14").Select()
